# Add 2 more private schools to the Racine schools directory (rows 52-53).
# Row 52 (Small World Montessori) was a mostly-blank row that already had
# the "A" column formatted (border + fill, no center alignment); row 53 was
# completely blank. Both need the same "center aligned, bordered" look as
# column A for the B/C/D grade-span columns, which we get by copying A52's
# number format (border + fill) onto each target cell and then centering it
# -- this reproduces the new style Excel created (cellXfs index 13) without
# hand-rolling borders cell by cell (which would create spurious new border
# entries instead of reusing the existing one).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A52").Copy()

# ---- Row 52: Small World Montessori (PK-KG) ----
$ws.Range("A52").Value = "Small World Montessori"

$ws.Range("B52").PasteSpecial(-4122)
$ws.Range("B52").Value = "PK"
$ws.Range("B52").HorizontalAlignment = -4108

$ws.Range("C52").Value = "KG"

$ws.Range("D52").PasteSpecial(-4122)
$ws.Range("D52").Value = 1
$ws.Range("D52").HorizontalAlignment = -4108

$ws.Range("I52").Value = "K9306402"

# ---- Row 53: St John's Lutheran School (PK-8) ----
$ws.Range("A53").PasteSpecial(-4122)
$ws.Range("A53").Value = "St John's Lutheran School"

$ws.Range("B53").PasteSpecial(-4122)
$ws.Range("B53").Value = "PK"
$ws.Range("B53").HorizontalAlignment = -4108

$ws.Range("C53").PasteSpecial(-4122)
$ws.Range("C53").Value = 8
$ws.Range("C53").HorizontalAlignment = -4108

$ws.Range("D53").PasteSpecial(-4122)
$ws.Range("D53").Value = 1
$ws.Range("D53").HorizontalAlignment = -4108

$ws.Range("E53").Value = 1

$ws.Range("I53").Value = 1511444

# Leave the selection where Excel would land after typing the last row.
$ws.Range("A53").Select()
